$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns to the right
$ws.Columns.Item(1).Insert()

# Set the new header in A1
$ws.Range("A1").Value = "DOCUMENTO"
